# Generate Report for Handback
# - Updates the "Status" text from "Ready for handoff" to "Handback transform failed"
#   for the 975cb303-... file (row 3) on both the Overview sheet and its per-language
#   status cells.
# - Records an "Error Detail" message explaining why the handback transform failed
#   on the zh-cn and de-de language sheets (column K, row 3 - the 975cb303 file row).

$wb = $excel.ActiveWorkbook

$newStatus = "Handback transform failed"

$zhDetail = "Handback file name: zqp5wjkf.rmr is different with handoff file name: 975cb303-2292-44de-a4f8-0a290b00db26.2196fd2d9a1d413eba96fd3e9134b5a311f8d0a7.zh-cn."
$deDetail = "Handback file name: zqp5wjkf.rmr is different with handoff file name: 975cb303-2292-44de-a4f8-0a290b00db26.2196fd2d9a1d413eba96fd3e9134b5a311f8d0a7.de-de."

# --- Overview sheet: update the status for the 975cb303 file (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $newStatus
$overview.Range("C3").Value = $newStatus

# --- zh-cn sheet: same Status text is shared for the 975cb303 file (row 3, column C) ---
# and add the Error Detail for that file (row 3, column K)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $newStatus
$zhcn.Range("K3").Value = $zhDetail

# --- de-de sheet: same Status text is shared for the 975cb303 file (row 3, column C) ---
# and add the Error Detail for that file (row 3, column K)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $newStatus
$dede.Range("K3").Value = $deDetail
